$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 20.10036642348715
$ws.Cells.Item(3, 2).Value = 19.53616503822142
$ws.Cells.Item(4, 2).Value = 19.18118461436417
$ws.Cells.Item(5, 2).Value = 19.0345523835764
$ws.Cells.Item(6, 2).Value = 19.01009037047504
$ws.Cells.Item(7, 2).Value = 19.17921483322056
$ws.Cells.Item(8, 2).Value = 19.90770846336206
$ws.Cells.Item(9, 2).Value = 21.2612471617243
$ws.Cells.Item(10, 2).Value = 22.20133159603876
$ws.Cells.Item(11, 2).Value = 22.61563291420337
$ws.Cells.Item(12, 2).Value = 22.77048476191025
$ws.Cells.Item(13, 2).Value = 22.73722687352303
$ws.Cells.Item(14, 2).Value = 22.62841407564251
$ws.Cells.Item(15, 2).Value = 22.5614948363401
$ws.Cells.Item(16, 2).Value = 22.17397728069621
$ws.Cells.Item(17, 2).Value = 21.93274304441507
$ws.Cells.Item(18, 2).Value = 21.79274127557897
$ws.Cells.Item(19, 2).Value = 21.74512807511792
$ws.Cells.Item(20, 2).Value = 21.95855310866448
$ws.Cells.Item(21, 2).Value = 22.66043112699017
$ws.Cells.Item(22, 2).Value = 23.107241519722
$ws.Cells.Item(23, 2).Value = 22.86989515375588
$ws.Cells.Item(24, 2).Value = 21.94688845935621
$ws.Cells.Item(25, 2).Value = 20.90405247994171
$ws.Cells.Item(2, 4).Value = 11.3164693344574
$ws.Cells.Item(3, 4).Value = 11.3669313530876
$ws.Cells.Item(4, 4).Value = 11.40007084026701
$ws.Cells.Item(5, 4).Value = 11.41411835830075
$ws.Cells.Item(6, 4).Value = 11.41648375390414
$ws.Cells.Item(7, 4).Value = 11.40025809064836
$ws.Cells.Item(8, 4).Value = 11.33342176000963
$ws.Cells.Item(9, 4).Value = 11.21942132384344
$ws.Cells.Item(10, 4).Value = 11.14601660603197
$ws.Cells.Item(11, 4).Value = 11.11485985390172
$ws.Cells.Item(12, 4).Value = 11.10338226342616
$ws.Cells.Item(13, 4).Value = 11.10583991068373
$ws.Cells.Item(14, 4).Value = 11.11390916010091
$ws.Cells.Item(15, 4).Value = 11.11889356623488
$ws.Cells.Item(16, 4).Value = 11.14809769901478
$ws.Cells.Item(17, 4).Value = 11.16658556166537
$ws.Cells.Item(18, 4).Value = 11.1774297378174
$ws.Cells.Item(19, 4).Value = 11.18113755421394
$ws.Cells.Item(20, 4).Value = 11.16459572149209
$ws.Cells.Item(21, 4).Value = 11.11153032478737
$ws.Cells.Item(22, 4).Value = 11.07871856971607
$ws.Cells.Item(23, 4).Value = 11.09605996904998
$ws.Cells.Item(24, 4).Value = 11.16549465765407
$ws.Cells.Item(25, 4).Value = 11.24843994939838
$ws.Cells.Item(2, 5).Value = 17.84191356875908
$ws.Cells.Item(3, 5).Value = 17.8087526950473
$ws.Cells.Item(4, 5).Value = 17.7921761536263
$ws.Cells.Item(5, 5).Value = 17.78638115811065
$ws.Cells.Item(6, 5).Value = 17.78547712363608
$ws.Cells.Item(7, 5).Value = 17.79209410233311
$ws.Cells.Item(8, 5).Value = 17.82969954253383
$ws.Cells.Item(9, 5).Value = 17.93306882784208
$ws.Cells.Item(10, 5).Value = 18.02649279174837
$ws.Cells.Item(11, 5).Value = 18.07265144794043
$ws.Cells.Item(12, 5).Value = 18.09064416742597
$ws.Cells.Item(13, 5).Value = 18.08674648822238
$ws.Cells.Item(14, 5).Value = 18.07412150519023
$ws.Cells.Item(15, 5).Value = 18.06645481314515
$ws.Cells.Item(16, 5).Value = 18.02354877434553
$ws.Cells.Item(17, 5).Value = 17.99815570075009
$ws.Cells.Item(18, 5).Value = 17.9838955430078
$ws.Cells.Item(19, 5).Value = 17.9791269791775
$ws.Cells.Item(20, 5).Value = 18.00082319389787
$ws.Cells.Item(21, 5).Value = 18.07781593916361
$ws.Cells.Item(22, 5).Value = 18.1311215141799
$ws.Cells.Item(23, 5).Value = 18.10240245404467
$ws.Cells.Item(24, 5).Value = 17.9996161643394
$ws.Cells.Item(25, 5).Value = 17.9019904851601
$ws.Cells.Item(2, 6).Value = 27.24940989735213
$ws.Cells.Item(3, 6).Value = 27.47001971590002
$ws.Cells.Item(4, 6).Value = 27.6149000655338
$ws.Cells.Item(5, 6).Value = 27.67630248125171
$ws.Cells.Item(6, 6).Value = 27.68664077981175
$ws.Cells.Item(7, 6).Value = 27.6157186041146
$ws.Cells.Item(8, 6).Value = 27.32351542856886
$ws.Cells.Item(9, 6).Value = 26.82563248149613
$ws.Cells.Item(10, 6).Value = 26.50611958709672
$ws.Cells.Item(11, 6).Value = 26.37092416045556
$ws.Cells.Item(12, 6).Value = 26.32119889471346
$ws.Cells.Item(13, 6).Value = 26.33184255929159
$ws.Cells.Item(14, 6).Value = 26.36680370514389
$ws.Cells.Item(15, 6).Value = 26.3884102075726
$ws.Cells.Item(16, 6).Value = 26.51516010195302
$ws.Cells.Item(17, 6).Value = 26.59552465970651
$ws.Cells.Item(18, 6).Value = 26.64270319663602
$ws.Cells.Item(19, 6).Value = 26.6588407446175
$ws.Cells.Item(20, 6).Value = 26.58687080546188
$ws.Cells.Item(21, 6).Value = 26.35649478106064
$ws.Cells.Item(22, 6).Value = 26.2145058612284
$ws.Cells.Item(23, 6).Value = 26.28949988130032
$ws.Cells.Item(24, 6).Value = 26.59078017386588
$ws.Cells.Item(25, 6).Value = 26.95223149841131
$ws.Cells.Item(2, 7).Value = 24.89018804700333
$ws.Cells.Item(3, 7).Value = 25.00564368637249
$ws.Cells.Item(4, 7).Value = 25.08958474275304
$ws.Cells.Item(5, 7).Value = 25.12704275728408
$ws.Cells.Item(6, 7).Value = 25.13345808542742
$ws.Cells.Item(7, 7).Value = 25.0900767914609
$ws.Cells.Item(8, 7).Value = 24.92727070138454
$ws.Cells.Item(9, 7).Value = 24.71283061053414
$ws.Cells.Item(10, 7).Value = 24.62077804489341
$ws.Cells.Item(11, 7).Value = 24.59341663728402
$ws.Cells.Item(12, 7).Value = 24.58516295589125
$ws.Cells.Item(13, 7).Value = 24.5868465067972
$ws.Cells.Item(14, 7).Value = 24.59269525750551
$ws.Cells.Item(15, 7).Value = 24.59655280290561
$ws.Cells.Item(16, 7).Value = 24.62286003986943
$ws.Cells.Item(17, 7).Value = 24.64273139673871
$ws.Cells.Item(18, 7).Value = 24.65552592567762
$ws.Cells.Item(19, 7).Value = 24.66009171508138
$ws.Cells.Item(20, 7).Value = 24.64047462261968
$ws.Cells.Item(21, 7).Value = 24.59091999184384
$ws.Cells.Item(22, 7).Value = 24.57082473793097
$ws.Cells.Item(23, 7).Value = 24.58041924510137
$ws.Cells.Item(24, 7).Value = 24.64149064360161
$ws.Cells.Item(25, 7).Value = 24.75944733760801
$ws.Cells.Item(2, 8).Value = 13.32108697707024
$ws.Cells.Item(3, 8).Value = 13.3883897808976
$ws.Cells.Item(4, 8).Value = 13.43274554819528
$ws.Cells.Item(5, 8).Value = 13.45158214885211
$ws.Cells.Item(6, 8).Value = 13.45475590039205
$ws.Cells.Item(7, 8).Value = 13.43299650423151
$ws.Cells.Item(8, 8).Value = 13.3436631474497
$ws.Cells.Item(9, 8).Value = 13.19258193542983
$ws.Cells.Item(10, 8).Value = 13.09634009074095
$ws.Cells.Item(11, 8).Value = 13.05577611827759
$ws.Cells.Item(12, 8).Value = 13.04087950524125
$ws.Cells.Item(13, 8).Value = 13.04406709310812
$ws.Cells.Item(14, 8).Value = 13.05454125539487
$ws.Cells.Item(15, 8).Value = 13.06101746431206
$ws.Cells.Item(16, 8).Value = 13.09905590149978
$ws.Cells.Item(17, 8).Value = 13.12321623598798
$ws.Cells.Item(18, 8).Value = 13.13741532251888
$ws.Cells.Item(19, 8).Value = 13.14227483344368
$ws.Cells.Item(20, 8).Value = 13.12061298795694
$ws.Cells.Item(21, 8).Value = 13.05145213291888
$ws.Cells.Item(22, 8).Value = 13.00895725902388
$ws.Cells.Item(23, 8).Value = 13.03138949696039
$ws.Cells.Item(24, 8).Value = 13.1217889539954
$ws.Cells.Item(25, 8).Value = 13.23086662253746
$ws.Cells.Item(2, 10).Value = 12.17260747219138
$ws.Cells.Item(3, 10).Value = 12.14449798796281
$ws.Cells.Item(4, 10).Value = 12.12924948587443
$ws.Cells.Item(5, 10).Value = 12.12354709207974
$ws.Cells.Item(6, 10).Value = 12.12263127108939
$ws.Cells.Item(7, 10).Value = 12.12917050292531
$ws.Cells.Item(8, 10).Value = 12.16250158725223
$ws.Cells.Item(9, 10).Value = 12.24354749174992
$ws.Cells.Item(10, 10).Value = 12.31226955302135
$ws.Cells.Item(11, 10).Value = 12.34543751176212
$ws.Cells.Item(12, 10).Value = 12.35826347250955
$ws.Cells.Item(13, 10).Value = 12.35548947494492
$ws.Cells.Item(14, 10).Value = 12.34648742810931
$ws.Cells.Item(15, 10).Value = 12.34100780717929
$ws.Cells.Item(16, 10).Value = 12.31013961556542
$ws.Cells.Item(17, 10).Value = 12.29168546641611
$ws.Cells.Item(18, 10).Value = 12.28125102771987
$ws.Cells.Item(19, 10).Value = 12.27774924620818
$ws.Cells.Item(20, 10).Value = 12.29363138090502
$ws.Cells.Item(21, 10).Value = 12.349124394905
$ws.Cells.Item(22, 10).Value = 12.3869386452361
$ws.Cells.Item(23, 10).Value = 12.3666177217372
$ws.Cells.Item(24, 10).Value = 12.29275108690959
$ws.Cells.Item(25, 10).Value = 12.21998165234376
$ws.Cells.Item(2, 12).Value = 12.66531499347237
$ws.Cells.Item(3, 12).Value = 12.15679027201557
$ws.Cells.Item(4, 12).Value = 11.83193337726028
$ws.Cells.Item(5, 12).Value = 11.69651615404923
$ws.Cells.Item(6, 12).Value = 11.67385129672434
$ws.Cells.Item(7, 12).Value = 11.83011918737934
$ws.Cells.Item(8, 12).Value = 12.4926724036742
$ws.Cells.Item(9, 12).Value = 13.68661553058802
$ws.Cells.Item(10, 12).Value = 14.49377191254199
$ws.Cells.Item(11, 12).Value = 14.84480841376603
$ws.Cells.Item(12, 12).Value = 14.97535146644205
$ws.Cells.Item(13, 12).Value = 14.94734377558039
$ws.Cells.Item(14, 12).Value = 14.8555964444193
$ws.Cells.Item(15, 12).Value = 14.79908594657228
$ws.Cells.Item(16, 12).Value = 14.47050017657939
$ws.Cells.Item(17, 12).Value = 14.26473889651674
$ws.Cells.Item(18, 12).Value = 14.14487420611554
$ws.Cells.Item(19, 12).Value = 14.10403183222347
$ws.Cells.Item(20, 12).Value = 14.28679996886698
$ws.Cells.Item(21, 12).Value = 14.88261010968844
$ws.Cells.Item(22, 12).Value = 15.25807062232969
$ws.Cells.Item(23, 12).Value = 15.05897447996225
$ws.Cells.Item(24, 12).Value = 14.27683104261209
$ws.Cells.Item(25, 12).Value = 13.37555565823421
$ws.Cells.Item(2, 15).Value = 19.758278515085
$ws.Cells.Item(3, 15).Value = 19.86999273839701
$ws.Cells.Item(4, 15).Value = 19.94491920226819
$ws.Cells.Item(5, 15).Value = 19.97703773349647
$ws.Cells.Item(6, 15).Value = 19.98246653316226
$ws.Cells.Item(7, 15).Value = 19.94534595442574
$ws.Cells.Item(8, 15).Value = 19.79547874004146
$ws.Cells.Item(9, 15).Value = 19.55216343726455
$ws.Cells.Item(10, 15).Value = 19.40466600649472
$ws.Cells.Item(11, 15).Value = 19.34444699671918
$ws.Cells.Item(12, 15).Value = 19.3226402689796
$ws.Cells.Item(13, 15).Value = 19.32729227701938
$ws.Cells.Item(14, 15).Value = 19.34263292025468
$ws.Cells.Item(15, 15).Value = 19.35215956376791
$ws.Cells.Item(16, 15).Value = 19.40874053948781
$ws.Cells.Item(17, 15).Value = 19.44521842292891
$ws.Cells.Item(18, 15).Value = 19.46684649401085
$ws.Cells.Item(19, 15).Value = 19.4742802633255
$ws.Cells.Item(20, 15).Value = 19.44126827539599
$ws.Cells.Item(21, 15).Value = 19.33809988399599
$ws.Cells.Item(22, 15).Value = 19.27648765090088
$ws.Cells.Item(23, 15).Value = 19.30883669748244
$ws.Cells.Item(24, 15).Value = 19.4430520928447
$ws.Cells.Item(25, 15).Value = 19.61252726398534
